# Update title slide text
$p = $ppt.ActivePresentation

# Slide 1: Title text change
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "US Software Comps – Oct 2025"

# Slide 7 (last slide): Title + content placeholder text changes
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Runs(1).Text = "Methodology & Sources"

$body = $s7.Shapes.Item(2).TextFrame.TextRange
$body.Paragraphs(1).Runs(1).Text = "Universe: 10 listed software comps; currency: USD; base: TTM."
$body.Paragraphs(2).Runs(1).Text = "Valuation: EV/EBITDA & P/E; medians and interquartile range (25–75th)."
$body.Paragraphs(3).Runs(1).Text = "Outliers reviewed; results illustrative. Sources: public filings & aggregators."
